# Auto-generated edit script updating cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.485.47"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.436.10"
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'572.15"
$ws.Range("E5").Value = "  -4.93%  "
$ws.Range("D6").Value = "'188.73"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").Value = "'0.607"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("D8").Value = "3.425.08"
$ws.Range("E8").Value = "  -5.19%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.201"
$ws.Range("E10").Value = "  -5.43%  "
$ws.Range("D11").Value = "'0.610"
$ws.Range("E11").Value = "  -5.55%  "
$ws.Range("D12").Value = "'50.65"
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("D13").Value = "'0.0000281"
$ws.Range("E13").Value = "  -7.42%  "
$ws.Range("D14").Value = "'8.99"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").Value = "3.973.92"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "'629.53"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "68.329.31"
$ws.Range("D18").Value = "3.450.91"
$ws.Range("E18").Value = "  -5.24%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "'12.13"
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").Value = "'18.00"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "'0.933"
$ws.Range("E22").Value = "  -6.54%  "
$ws.Range("D23").Value = "'17.63"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "'5.31"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'98.28"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("E26").Value = "  -8.53%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.81"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'6.06"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "'9.77"
$ws.Range("E29").Value = "  -7.92%  "
$ws.Range("D30").Value = "'9.11"
$ws.Range("E30").Value = "  -6.04%  "
$ws.Range("D31").Value = "'32.05"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").Value = "'4.13"
$ws.Range("E32").Value = "  -12.12%  "
$ws.Range("D33").Value = "'6.62"
$ws.Range("E33").Value = "  -8.71%  "
$ws.Range("D34").Value = "'11.47"
$ws.Range("E34").Value = "  -6.55%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = "  -7.49%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'60.42"
$ws.Range("E36").Value = "  -4.46%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "3.620.15"
$ws.Range("E38").Value = "  -8.14%  "
$ws.Range("D39").Value = "0.0₃0772"
$ws.Range("E39").Value = "  -13.03%  "
$ws.Range("D40").Value = "'495.97"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").Value = "'2.85"
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("D43").Value = "'0.364"
$ws.Range("E43").Value = "  -6.23%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "'3.49"
$ws.Range("E44").Value = "  +69.54%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.132"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").Value = "'33.93"
$ws.Range("E46").Value = "  -7.52%  "
$ws.Range("D47").Value = "'0.0433"
$ws.Range("E47").Value = "  -5.82%  "
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("D49").Value = "'2.78"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "  -0.51%  "
